$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "64.141.89"
$ws.Range("E2").Value = "  +0.87%  "
Set-TextValue $ws.Range("D3") "3.136.88"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  +0.28%  "
Set-TextValue $ws.Range("D5") "601.08"
$ws.Range("E5").Value = "  -1.02%  "
Set-TextValue $ws.Range("D6") "143.73"
$ws.Range("E6").Value = "  -1.32%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.28%  "
Set-TextValue $ws.Range("D8") "3.127.26"
$ws.Range("E8").Value = "  +0.28%  "
Set-TextValue $ws.Range("D9") "0.521"
$ws.Range("E9").Value = "  +0.01%  "
Set-TextValue $ws.Range("D10") "0.149"
$ws.Range("E10").Value = "  -0.79%  "
Set-TextValue $ws.Range("D11") "5.35"
$ws.Range("E11").Value = "  +0.95%  "
Set-TextValue $ws.Range("D12") "0.468"
$ws.Range("E12").Value = "  -0.44%  "
Set-TextValue $ws.Range("D13") "0.0000253"
$ws.Range("E13").Value = "  +0.48%  "
Set-TextValue $ws.Range("D14") "35.33"
$ws.Range("E14").Value = "  +0.04%  "
Set-TextValue $ws.Range("D15") "3.671.20"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("E16").Value = "  +2.40%  "
Set-TextValue $ws.Range("D17") "64.318.34"
$ws.Range("E17").Value = "  +0.91%  "
Set-TextValue $ws.Range("D18") "3.157.03"
$ws.Range("E18").Value = "  +0.85%  "
Set-TextValue $ws.Range("D19") "6.85"
$ws.Range("E19").Value = "  +0.31%  "
Set-TextValue $ws.Range("D20") "480.68"
$ws.Range("E20").Value = "  +1.15%  "
Set-TextValue $ws.Range("D21") "14.62"
$ws.Range("E21").Value = "  +0.50%  "
Set-TextValue $ws.Range("D22") "0.710"
$ws.Range("E22").Value = "  -0.04%  "
Set-TextValue $ws.Range("D23") "7.65"
$ws.Range("E23").Value = "  -2.44%  "
Set-TextValue $ws.Range("D24") "87.97"
$ws.Range("E24").Value = "  +5.44%  "
Set-TextValue $ws.Range("D25") "13.40"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("E26").Value = "  -0.04%  "
Set-TextValue $ws.Range("D27") "2.75"
$ws.Range("E27").Value = "  -1.53%  "
Set-TextValue $ws.Range("D28") "8.35"
$ws.Range("E28").Value = "  -1.45%  "
Set-TextValue $ws.Range("D29") "7.12"
$ws.Range("E29").Value = "  +1.48%  "
Set-TextValue $ws.Range("D30") "2.08"
$ws.Range("E30").Value = "  +0.45%  "
Set-TextValue $ws.Range("D31") "0.111"
$ws.Range("E31").Value = "  -8.14%  "
Set-TextValue $ws.Range("D34") "2.67"
$ws.Range("E34").Value = "  -0.62%  "
Set-TextValue $ws.Range("D35") "1.10"
$ws.Range("E35").Value = "  -1.60%  "
Set-TextValue $ws.Range("D36") "6.03"
$ws.Range("E36").Value = "  +1.52%  "
Set-TextValue $ws.Range("D37") "0.0₃0756"
$ws.Range("E37").Value = "  -3.97%  "
Set-TextValue $ws.Range("D38") "52.68"
$ws.Range("E38").Value = "  -0.34%  "
Set-TextValue $ws.Range("D39") "3.01"
$ws.Range("E39").Value = "  -0.04%  "
Set-TextValue $ws.Range("D40") "440.05"
$ws.Range("E40").Value = "  -3.55%  "
Set-TextValue $ws.Range("D41") "0.0394"
$ws.Range("E41").Value = "  +0.26%  "
Set-TextValue $ws.Range("D42") "0.118"
$ws.Range("E42").Value = "  +0.14%  "
Set-TextValue $ws.Range("D43") "8.26"
$ws.Range("E43").Value = "  -0.66%  "
Set-TextValue $ws.Range("D44") "2.869.18"
$ws.Range("E44").Value = "  +0.66%  "
Set-TextValue $ws.Range("D45") "0.262"
$ws.Range("E45").Value = "  -1.87%  "
Set-TextValue $ws.Range("D46") "2.46"
$ws.Range("E46").Value = "  +1.80%  "
Set-TextValue $ws.Range("D47") "2.22"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("E48").Value = "  -0.01%  "
Set-TextValue $ws.Range("D49") "25.92"
$ws.Range("E49").Value = "  -1.21%  "
Set-TextValue $ws.Range("D50") "0.113"
$ws.Range("E50").Value = "  +0.21%  "
Set-TextValue $ws.Range("D51") "121.50"
$ws.Range("E51").Value = "  +2.17%  "

# Rows 32 and 33: coin order/content swap (EthereumClassic <-> FirstDigitalUSD) with updated figures
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D32") "1.00"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D33") "27.01"
$ws.Range("E33").Value = "  +3.20%  "
